$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 96 (weekly update: pushes existing rows 96:187 -> 97:188)
$ws.Rows.Item(96).EntireRow.Insert()

# Populate the new row 96 - duplicate of (old) row 96 data but with the new, more
# recent reporting date (2022-01-28 -> serial 44589)
$ws.Cells.Item(96, 1).Value = 7
$ws.Cells.Item(96, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(96, 3).Value = "Ñuble"
$ws.Cells.Item(96, 4).Value = 44589
$ws.Cells.Item(96, 5).Value = 16
$ws.Cells.Item(96, 6).Value = 100112003
$ws.Cells.Item(96, 7).Value = "Ajo"
$ws.Cells.Item(96, 8).Value = "Chino"
$ws.Cells.Item(96, 9).Value = "Primera"
$ws.Cells.Item(96, 10).Value = 100
$ws.Cells.Item(96, 11).Value = 18000
$ws.Cells.Item(96, 12).Value = 19000
$ws.Cells.Item(96, 13).Value = 18500
$ws.Cells.Item(96, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(96, 15).Value = "China"
$ws.Cells.Item(96, 16).Value = 1850
$ws.Cells.Item(96, 17).Value = 10
$ws.Cells.Item(96, 18).Value = "Hortaliza"
